$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metrics")
$src = $ws.Range("E2")
$dst = $ws.Range("E28")
$dst.Value2 = "number_of_connections_to_workstations_for_server"
$src.Copy()
$dst.PasteSpecial(-4122)  # xlPasteFormats
